# TimeLog_ConnorPeper.xlsx - "Lots of new documents"
#
# Updates week-3 row (row 9) with hours worked and an activity note, appends
# a short remark to the week-2 activity note, records a new shared-string
# for the meeting note, and moves the selection / view position to reflect
# where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 2 (row 8) activity note gains a trailing remark.
$ws.Range("F8").Value = "Worked on the Vision document with the rest of the team. Made some demo ER-Diagrams and UML diagrams in spare time. Attended team meetings"

# Week 3 (row 9): hours worked this week, plus the new activity note.
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = "Meeting with professor. Met with team to choose SCRUM master and Product Owner."

# Reposition the view / selection to where the user left off editing.
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 5
$win.ScrollColumn = 4
$ws.Range("E10").Select()
